# Bug fixes based on Vidhya's feedback:
# - 24-hour day is now used (instead of 12 hour day)
# - Dates previously could be accidently saved with invalid byte order
#
# Update the orderDate (column C) and orderTime (column D) sample values,
# and move the active selection to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# orderDate column (C) corrections
$ws.Range("C2").Value = 367
$ws.Range("C3").Value = 764
$ws.Range("C4").Value = 37683
$ws.Range("C6").Value = 239360

# orderTime column (D) corrections (24-hour values)
$ws.Range("D3").Value = 0.99930555555555556
$ws.Range("D5").Value = 0.78055555555555556

# Move the active selection
$ws.Range("C7").Select()
